$wb = $excel.ActiveWorkbook

# --- Update data values on sheet "Grade 1 or 2" (before rename) ---
$ws = $wb.Worksheets.Item("Grade 1 or 2")
$ws.Range("C8").Value = -0.11931818181818188
$ws.Range("D8").Value = -0.06896551724137924
$ws.Range("E8").Value = -0.02316602316602316
$ws.Range("F8").Value = 0.02290076335877862
$ws.Range("C10").Value = -0.09006622516556292
$ws.Range("D10").Value = -0.06389214536928484
$ws.Range("E10").Value = -0.022435897435897505
$ws.Range("F10").Value = 0.030519790176442484
$ws.Range("C11").Value = -0.011054024599096988
$ws.Range("D11").Value = 0.008026853473438338
$ws.Range("E11").Value = 0.01972611672644278
$ws.Range("F11").Value = 0.03985615822595139
$ws.Range("C13").Value = -0.021523178807947053
$ws.Range("D13").Value = -0.02505297270347751
$ws.Range("E13").Value = -0.011617846990483346
$ws.Range("F13").Value = -0.01025791324736226
$ws.Range("C14").Value = 0.07799290973547847
$ws.Range("D14").Value = 0.0510913888222595
$ws.Range("E14").Value = 0.026333407721490686
$ws.Range("F14").Value = -0.027879835746704176
$ws.Range("C15").Value = 0.10100130605137148
$ws.Range("D15").Value = 0.15255423372988114
$ws.Range("E15").Value = 0.35986653956148723
$ws.Range("F15").Value = 0.28040540540540543
$ws.Range("C16").Value = 0.021827593353654483
$ws.Range("D16").Value = -0.00782888298615962
$ws.Range("E16").Value = -0.0015085012222495198
$ws.Range("F16").Value = -0.02590303985402069
$ws.Range("C17").Value = -0.040485829959514205
$ws.Range("D17").Value = 0.006276150627614952
$ws.Range("E17").Value = 0.008510638297872289
$ws.Range("F17").Value = 0.03703703703703707
$ws.Range("C18").Value = -0.010249839846252341
$ws.Range("D18").Value = 0.008721359940872161
$ws.Range("E18").Value = 0.015972812234494416
$ws.Range("F18").Value = 0.035277947464874615
$ws.Range("C19").Value = -0.2173913043478261
$ws.Range("D19").Value = 0.11538461538461536
$ws.Range("E19").Value = -0.09090909090909094
$ws.Range("F19").Value = 0.5000000000000001
$ws.Range("C20").Value = -0.01185822256971592
$ws.Range("D20").Value = 0.01580611169652267
$ws.Range("E20").Value = 0.0013517166801838346
$ws.Range("F20").Value = 0.04926108374384229

# --- Update data values on sheet "Home Owner" (before rename) ---
$ws = $wb.Worksheets.Item("Home Owner")
$ws.Range("C2").Value = 0.02087332053742799
$ws.Range("D2").Value = -0.022929936305732503
$ws.Range("E2").Value = 0.013556388808768336
$ws.Range("F2").Value = -0.06412175648702595
$ws.Range("C3").Value = 0.06816030293467973
$ws.Range("D3").Value = 0.0485232067510549
$ws.Range("E3").Value = 0.08266414737836555
$ws.Range("F3").Value = -0.04821638573108587
$ws.Range("C4").Value = 0.04176771759633526
$ws.Range("D4").Value = -0.0015455950540958282
$ws.Range("E4").Value = 0.03672612801678912
$ws.Range("F4").Value = -0.059173357121617524
$ws.Range("C5").Value = 0.028028028028027997
$ws.Range("D5").Value = -0.01590348231423095
$ws.Range("E5").Value = 0.01668653158522048
$ws.Range("F5").Value = -0.06338397066526978
$ws.Range("C6").Value = -0.01378070701018572
$ws.Range("D6").Value = -0.06324110671936757
$ws.Range("E6").Value = 0.010928961748633977
$ws.Range("F6").Value = -0.14440639269406386
$ws.Range("C7").Value = 0.00509554140127393
$ws.Range("D7").Value = -0.02058590657165472
$ws.Range("E7").Value = 0.0071820870299958055
$ws.Range("F7").Value = -0.06416217221498972
$ws.Range("C8").Value = 0.0058823529411765035
$ws.Range("D8").Value = -0.04481132075471696
$ws.Range("E8").Value = 0.018817204301075436
$ws.Range("F8").Value = -0.1379310344827585
$ws.Range("C9").Value = 0.0032258064516128837
$ws.Range("D9").Value = -0.05346350534635058
$ws.Range("E9").Value = 0.005485463521667586
$ws.Range("F9").Value = -0.13836734693877545
$ws.Range("C10").Value = -0.004013220018885693
$ws.Range("D10").Value = -0.05597897503285147
$ws.Range("E10").Value = -0.007849293563579284
$ws.Range("F10").Value = -0.13587715216379712
$ws.Range("C11").Value = 0.0282478686855834
$ws.Range("D11").Value = 0.026081888387002806
$ws.Range("E11").Value = 0.0017231476163124338
$ws.Range("F11").Value = 0.019201462968607048
$ws.Range("C12").Value = -0.05613648871766641
$ws.Range("D12").Value = -0.0762112139357649
$ws.Range("E12").Value = 0.01731879409878114
$ws.Range("F12").Value = -0.14517625231910955
$ws.Range("C13").Value = -0.029202739372011968
$ws.Range("D13").Value = -0.03418367346938781
$ws.Range("E13").Value = -0.019283403831798905
$ws.Range("F13").Value = -0.025254762959680998
$ws.Range("C14").Value = -0.29859719438877746
$ws.Range("D14").Value = -0.18413021363173962
$ws.Range("E14").Value = -0.08339889850511396
$ws.Range("F14").Value = -0.1287410926365795
$ws.Range("C15").Value = 0.27280064568200174
$ws.Range("D15").Value = 0.34804490902051893
$ws.Range("E15").Value = 0.22982975573649145
$ws.Range("F15").Value = 0.1227085671530116
$ws.Range("C16").Value = 0.02671969863946943
$ws.Range("D16").Value = 0.008530350535265227
$ws.Range("E16").Value = 0.01378155823002446
$ws.Range("F16").Value = -0.012935748883643993
$ws.Range("C17").Value = -0.0815450643776825
$ws.Range("D17").Value = -0.03389830508474589
$ws.Range("E17").Value = -0.03984063745019924
$ws.Range("F17").Value = -0.007462686567164263
$ws.Range("C18").Value = 0.018231540565177774
$ws.Range("D18").Value = 0.021619927933573608
$ws.Range("E18").Value = -0.004089669797031257
$ws.Range("F18").Value = 0.029387485553904635
$ws.Range("C19").Value = -0.2000000000000001
$ws.Range("D19").Value = 0.03846153846153856
$ws.Range("E19").Value = -0.14285714285714293
$ws.Range("F19").Value = 0.8333333333333336
$ws.Range("C20").Value = -0.013757861635220065
$ws.Range("D20").Value = -0.015608521408985483
$ws.Range("E20").Value = -0.007016596950478994
$ws.Range("F20").Value = 0.026975763962065403

# --- Update data values on sheet "Inactive" (before rename) ---
$ws = $wb.Worksheets.Item("Inactive")
$ws.Range("C2").Value = -0.041371158392434895
$ws.Range("D2").Value = -0.003875968992247994
$ws.Range("E2").Value = -0.04773101340230422
$ws.Range("F2").Value = -0.05187376725838262
$ws.Range("C3").Value = 0.025769506084466858
$ws.Range("D3").Value = 0.07115531752104057
$ws.Range("E3").Value = -0.03989361702127663
$ws.Range("F3").Value = -0.019656019656019708
$ws.Range("C4").Value = -0.016727069974909407
$ws.Range("D4").Value = 0.02306248323947438
$ws.Range("E4").Value = -0.04545454545454534
$ws.Range("F4").Value = -0.04087385482734326
$ws.Range("C5").Value = -0.033558727773603726
$ws.Range("D5").Value = 0.0035104142288790107
$ws.Range("E5").Value = -0.04779233381853468
$ws.Range("F5").Value = -0.04935711323102456
$ws.Range("C6").Value = -0.10015649452269182
$ws.Range("D6").Value = -0.02412380518889394
$ws.Range("E6").Value = -0.09493670886075943
$ws.Range("F6").Value = -0.12795698924731194
$ws.Range("C7").Value = -0.05476107429368675
$ws.Range("D7").Value = -0.009062448508815208
$ws.Range("E7").Value = -0.05043782837127835
$ws.Range("F7").Value = -0.05821359798548363
$ws.Range("C8").Value = -0.09090909090909088
$ws.Range("D8").Value = -0.0060422960725075164
$ws.Range("E8").Value = -0.11418685121107258
$ws.Range("F8").Value = -0.13492063492063486
$ws.Range("C9").Value = -0.13152173913043486
$ws.Range("D9").Value = -0.06272000000000003
$ws.Range("E9").Value = -0.12830735773831095
$ws.Range("F9").Value = -0.14694704839118317
$ws.Range("C10").Value = -0.17076923076923076
$ws.Range("D10").Value = -0.11434591074506806
$ws.Range("E10").Value = -0.14221652786675207
$ws.Range("F10").Value = -0.1553690212226797
$ws.Range("C11").Value = -0.027893639207507787
$ws.Range("D11").Value = -0.03948080043266631
$ws.Range("E11").Value = -0.06928345626975761
$ws.Range("F11").Value = -0.04184704184704176
$ws.Range("C12").Value = -0.04272727272727269
$ws.Range("D12").Value = 0.03519749706687528
$ws.Range("E12").Value = -0.028864059590316498
$ws.Range("F12").Value = -0.08900687070580877
$ws.Range("C13").Value = -0.04873044370351376
$ws.Range("D13").Value = -0.05791067373202112
$ws.Range("E13").Value = -0.029418975239029203
$ws.Range("F13").Value = -0.03518761863045686
$ws.Range("C14").Value = -0.1246458923512747
$ws.Range("D14").Value = -0.13398553806890684
$ws.Range("E14").Value = -0.05753514220333444
$ws.Range("F14").Value = -0.056494447126991865
$ws.Range("C15").Value = 0.16017885323513945
$ws.Range("D15").Value = 0.2575849514563106
$ws.Range("E15").Value = -0.01649132783622413
$ws.Range("F15").Value = 0.008164275111330948
$ws.Range("C16").Value = 0.05687652054281785
$ws.Range("D16").Value = 0.02891586899845863
$ws.Range("E16").Value = 0.04448540074282669
$ws.Range("F16").Value = 0.00847077988189786
$ws.Range("C17").Value = -0.1778846153846153
$ws.Range("D17").Value = -0.10294117647058829
$ws.Range("E17").Value = -0.09090909090909091
$ws.Range("F17").Value = -0.05263157894736837
$ws.Range("C18").Value = -0.038888888888888924
$ws.Range("D18").Value = -0.07075038284839207
$ws.Range("E18").Value = -0.08588252190147935
$ws.Range("F18").Value = -0.04096889952153107
$ws.Range("C19").Value = -0.2500000000000001
$ws.Range("D19").Value = -0.47058823529411764
$ws.Range("E19").Value = -0.2857142857142857
$ws.Range("F19").Value = 0.39999999999999997
$ws.Range("C20").Value = -0.02927019643554062
$ws.Range("D20").Value = -0.06890720777848239
$ws.Range("E20").Value = -0.023431294678316044
$ws.Range("F20").Value = -0.014683975313896492

# --- Update data values on sheet "LLTI" (before rename) ---
$ws = $wb.Worksheets.Item("LLTI")
$ws.Range("C2").Value = -0.22268907563025206
$ws.Range("D2").Value = -0.1268260292164675
$ws.Range("E2").Value = -0.057645631067961216
$ws.Range("F2").Value = -0.28811121764141906
$ws.Range("C3").Value = -0.2547085201793723
$ws.Range("D3").Value = -0.10728476821192048
$ws.Range("E3").Value = 0.03111111111111108
$ws.Range("F3").Value = -0.32409739714525604
$ws.Range("C4").Value = -0.23458540042523035
$ws.Range("D4").Value = -0.11830742659758207
$ws.Range("E4").Value = -0.02841781874039945
$ws.Range("F4").Value = -0.30232558139534876
$ws.Range("C5").Value = -0.22539682539682543
$ws.Range("D5").Value = -0.12243401759530786
$ws.Range("E5").Value = -0.05296343001261029
$ws.Range("F5").Value = -0.292197858235594
$ws.Range("C6").Value = -0.32423208191126274
$ws.Range("D6").Value = -0.13513513513513511
$ws.Range("E6").Value = -0.019543973941368184
$ws.Range("F6").Value = -0.4437229437229437
$ws.Range("C7").Value = -0.1961038961038961
$ws.Range("D7").Value = -0.0957297043641483
$ws.Range("E7").Value = -0.033949290932531115
$ws.Range("F7").Value = -0.2613833854709768
$ws.Range("C8").Value = -0.3295454545454546
$ws.Range("D8").Value = -0.1447368421052632
$ws.Range("E8").Value = -0.044943820224719024
$ws.Range("F8").Value = -0.4318181818181818
$ws.Range("C9").Value = -0.3130841121495327
$ws.Range("D9").Value = -0.13043478260869562
$ws.Range("E9").Value = -0.02272727272727275
$ws.Range("F9").Value = -0.4303405572755418
$ws.Range("C10").Value = -0.29255319148936176
$ws.Range("D10").Value = -0.1149068322981366
$ws.Range("E10").Value = 0.006345177664974625
$ws.Range("F10").Value = -0.4232804232804233
$ws.Range("C11").Value = -0.1515723270440252
$ws.Range("D11").Value = -0.08507386653082014
$ws.Range("E11").Value = -0.0458874458874458
$ws.Range("F11").Value = -0.1660377358490565
$ws.Range("C12").Value = -0.3485714285714286
$ws.Range("D12").Value = -0.155223880597015
$ws.Range("E12").Value = -0.018229166666666647
$ws.Range("F12").Value = -0.4661157024793388
$ws.Range("C13").Value = -0.028820210939416193
$ws.Range("D13").Value = -0.024024390243902367
$ws.Range("E13").Value = -0.009549461312438686
$ws.Range("F13").Value = -0.012416979497545395
$ws.Range("C14").Value = 0.0452240646841168
$ws.Range("D14").Value = 0.03366906474820157
$ws.Range("E14").Value = 0.02627422828427851
$ws.Range("F14").Value = 0.0788629069234296
$ws.Range("C15").Value = -0.14483394833948332
$ws.Range("D15").Value = 0.10987261146496821
$ws.Range("E15").Value = 0.0957642725598526
$ws.Range("F15").Value = -0.12362204724409454
$ws.Range("C16").Value = 0.042091301354913374
$ws.Range("D16").Value = 0.00973901113178328
$ws.Range("E16").Value = -0.017420116464787687
$ws.Range("F16").Value = 0.0016939840673932144
$ws.Range("C17").Value = -0.09226190476190473
$ws.Range("D17").Value = -0.0017953321364452938
$ws.Range("E17").Value = 0.1010452961672475
$ws.Range("F17").Value = 0.01257861635220118
$ws.Range("C18").Value = -0.15659881812212728
$ws.Range("D18").Value = -0.09379093528949436
$ws.Range("E18").Value = -0.03988868274582567
$ws.Range("F18").Value = -0.15137956748695006
$ws.Range("C19").Value = -0.29629629629629634
$ws.Range("D19").Value = 0.233009708737864
$ws.Range("E19").Value = 0.5833333333333335
$ws.Range("F19").Value = 1.3333333333333333
$ws.Range("C20").Value = -0.013063685466650014
$ws.Range("D20").Value = -0.0179442508710801
$ws.Range("E20").Value = 0.007094943240454011
$ws.Range("F20").Value = -0.005648720211826974

# --- Update data values on sheet "No Car" (before rename) ---
$ws = $wb.Worksheets.Item("No Car")
$ws.Range("C2").Value = -0.12386706948640488
$ws.Range("D2").Value = -0.09311740890688261
$ws.Range("E2").Value = -0.11168032786885253
$ws.Range("F2").Value = -0.11807120324470471
$ws.Range("C3").Value = -0.1378839590443686
$ws.Range("D3").Value = -0.05607476635514011
$ws.Range("E3").Value = -0.10044642857142856
$ws.Range("F3").Value = -0.11163416274377945
$ws.Range("C4").Value = -0.12794033275960995
$ws.Range("D4").Value = -0.08327501749475154
$ws.Range("E4").Value = -0.11190903650508674
$ws.Range("F4").Value = -0.11695906432748535
$ws.Range("C5").Value = -0.1229551451187336
$ws.Range("D5").Value = -0.08785046728971964
$ws.Range("E5").Value = -0.11152025249868489
$ws.Range("F5").Value = -0.11856400566839873
$ws.Range("C6").Value = -0.24791086350974934
$ws.Range("D6").Value = -0.1409395973154362
$ws.Range("E6").Value = -0.17971014492753634
$ws.Range("F6").Value = -0.22601279317697226
$ws.Range("C7").Value = -0.120923415170392
$ws.Range("D7").Value = -0.07603305785123968
$ws.Range("E7").Value = -0.0945945945945947
$ws.Range("F7").Value = -0.11331351172778327
$ws.Range("C8").Value = -0.2540983606557377
$ws.Range("D8").Value = -0.11578947368421055
$ws.Range("E8").Value = -0.17391304347826086
$ws.Range("F8").Value = -0.2251655629139073
$ws.Range("C9").Value = -0.24190800681431013
$ws.Range("D9").Value = -0.12147505422993493
$ws.Range("E9").Value = -0.16967509025270755
$ws.Range("F9").Value = -0.2169680111265647
$ws.Range("C10").Value = -0.23205506391347094
$ws.Range("D10").Value = -0.13135068153655513
$ws.Range("E10").Value = -0.1671891327063741
$ws.Range("F10").Value = -0.20734693877551014
$ws.Range("C11").Value = -0.08017334777898148
$ws.Range("D11").Value = -0.057191634656423365
$ws.Range("E11").Value = -0.060800790904597155
$ws.Range("F11").Value = -0.0910458991723099
$ws.Range("C12").Value = -0.2712933753943218
$ws.Range("D12").Value = -0.17589576547231275
$ws.Range("E12").Value = -0.20312500000000006
$ws.Range("F12").Value = -0.2505050505050505
$ws.Range("C13").Value = -0.03513022410660209
$ws.Range("D13").Value = -0.019123020706455465
$ws.Range("E13").Value = -0.01403808593750008
$ws.Range("F13").Value = -0.011384925781812808
$ws.Range("C14").Value = -0.015508395522388033
$ws.Range("D14").Value = 0.020378632581695304
$ws.Range("E14").Value = -0.0017736786094360373
$ws.Range("F14").Value = 0.032736177342441004
$ws.Range("C15").Value = -0.188715953307393
$ws.Range("D15").Value = -0.10402684563758384
$ws.Range("E15").Value = 0.019438444924406006
$ws.Range("F15").Value = -0.13351498637602183
$ws.Range("C16").Value = 0.07682308232704019
$ws.Range("D16").Value = 0.030466301953524833
$ws.Range("E16").Value = 0.022615272949535734
$ws.Range("F16").Value = 0.018483885610531072
$ws.Range("C17").Value = -0.182089552238806
$ws.Range("D17").Value = -0.06989247311827958
$ws.Range("E17").Value = -0.03952569169960478
$ws.Range("F17").Value = -0.044303797468354486
$ws.Range("C18").Value = -0.056547619047619097
$ws.Range("D18").Value = -0.03595812471552125
$ws.Range("E18").Value = -0.03638368246968023
$ws.Range("F18").Value = -0.06636060100166943
$ws.Range("C19").Value = -0.41379310344827586
$ws.Range("D19").Value = -0.2000000000000001
$ws.Range("E19").Value = -0.11111111111111105
$ws.Range("F19").Value = -0.2222222222222222
$ws.Range("C20").Value = -0.02832574607991904
$ws.Range("D20").Value = -0.0482013113592061
$ws.Range("E20").Value = -0.010800842992623795
$ws.Range("F20").Value = -0.025187566988210078

# --- Rename sheets (order matters to avoid name collisions) ---
$wb.Worksheets.Item("Grade 1 or 2").Name = "Grades 3 to 5"
$wb.Worksheets.Item("Home Owner").Name = "Has Car"
$wb.Worksheets.Item("Inactive").Name = "Home Owner"
$wb.Worksheets.Item("LLTI").Name = "Inactive"
$wb.Worksheets.Item("No Car").Name = "LLTI"